# Prologue.xlsx: point the scene-load call at the real start scene instead
# of the placeholder map, and leave the cursor on the cell that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 'Wait(1f,LoadSceneByEnum("GameStartScene"));'

[void]$ws.Range("C6").Select()
